# Rename the "Transportadora" header column to "Importadora".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reserva")

# The header row (row 1) holds column titles in A1:H1; the third column
# (C1) currently reads "Transportadora" and should read "Importadora".
$ws.Range("C1").Value = "Importadora"
